$wb = $excel.ActiveWorkbook

# Sheet "展览" updates (column F = "想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F11").Value = 2087
$ws1.Range("F13").Value = 1433
$ws1.Range("F26").Value = 2935
$ws1.Range("F28").Value = 3289
$ws1.Range("F35").Value = 382
$ws1.Range("F36").Value = 244

# Sheet "全部类型" updates (column F = "想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F13").Value = 2087
$ws4.Range("F15").Value = 1434
$ws4.Range("F19").Value = 4036
$ws4.Range("F30").Value = 2937
$ws4.Range("F32").Value = 3289
$ws4.Range("F39").Value = 382
$ws4.Range("F40").Value = 244
